# Russian-DOOM raven_charmap.xlsx: add a "Hexen" sheet (copy of "Heretic")
# with a couple of corrected glyph-map cells, fix the "02" row on both
# sheets to point at the backtick placeholder instead of a stray quote,
# and update the selection shown on the "Heretic" tab.

$wb = $excel.ActiveWorkbook
$heretic = $wb.Worksheets.Item("Heretic")

# --- Fix row 5 ("02") on Heretic: it should show the unused/backtick
#     placeholder (like the other "free" rows), not a literal quote char.
$heretic.Range("C10").Copy() | Out-Null
$heretic.Range("C5").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$heretic.Range("C5").Value = "``"

# --- Duplicate the sheet to create "Hexen" right after "Heretic".
$heretic.Copy($null, $heretic) | Out-Null
$hexen = $wb.Worksheets.Item(2)
$hexen.Name = "Hexen"

# --- Hexen-specific glyph fixes (rows 7 and 9 become "free" slots too).
$hexen.Range("C10").Copy() | Out-Null
$hexen.Range("C7").PasteSpecial(-4122) | Out-Null
$hexen.Range("C7").Value = "liber oscura"

$hexen.Range("C10").Copy() | Out-Null
$hexen.Range("C9").PasteSpecial(-4122) | Out-Null
$hexen.Range("C9").Value = "``"

$hexen.Columns("C").ColumnWidth = 10.4518229166

# --- Restore "Heretic" as the active/selected tab, with C5 selected.
$heretic.Activate()
$heretic.Range("C5").Select() | Out-Null
